# Fix test fixture data: rename header and correct the "123" mobile-number
# test case result (login validation was incorrectly failing it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header "mobile number" -> "mobile_number"
$ws.Range("A1").Value = "mobile_number"

# Row 4 ("123") should report "pass", not "fail"
$ws.Range("B4").Value = "pass"

# Move the active selection to the cell that was corrected
[void]$ws.Range("B4").Select()
